{"js": "// Map of old text -> new text for this edit. Every old value is unique\n// within the document, so a direct search-and-replace per pair is safe.\nconst replacements = [\n  [\"2024-12-11 Wednesday\", \"2024-12-12 Thursday\"],\n  [\"111\u00d76=666\", \"840\u00d79=7560\"],\n  [\"534\u00d77=3738\", \"832\u00d79=7488\"],\n  [\"468\u00d73=1404\", \"223\u00d74=892\"],\n  [\"390\u00d72=780\", \"858\u00d79=7722\"],\n  [\"891\u00d78=7128\", \"105\u00d79=945\"],\n  [\"292\u00d79=2628\", \"523\u00d76=3138\"],\n  [\"973\u00d75=4865\", \"528\u00d79=4752\"],\n  [\"379\u00d72=758\", \"860\u00d75=4300\"],\n  [\"380\u00d77=2660\", \"921\u00d73=2763\"],\n  [\"704\u00d77=4928\", \"485\u00d76=2910\"],\n  [\"470\u00d74=1880\", \"383\u00d76=2298\"],\n  [\"545\u00d75=2725\", \"144\u00d74=576\"],\n  [\"755\u00d76=4530\", \"419\u00d76=2514\"],\n  [\"389\u00d79=3501\", \"164\u00d75=820\"],\n  [\"801\u00d78=6408\", \"401\u00d74=1604\"],\n  [\"337\u00d77=2359\", \"975\u00d74=3900\"],\n  [\"345\u00d76=2070\", \"248\u00d72=496\"],\n  [\"486\u00d78=3888\", \"491\u00d74=1964\"],\n  [\"891\u00d79=8019\", \"125\u00d73=375\"],\n  [\"446\u00d77=3122\", \"703\u00d76=4218\"],\n  [\"696\u00d74=2784\", \"242\u00d73=726\"],\n  [\"205\u00d78=1640\", \"536\u00d77=3752\"],\n  [\"313\u00d72=626\", \"456\u00d74=1824\"],\n  [\"380\u00d73=1140\", \"305\u00d73=915\"],\n  [\"287\u00d72=574\", \"149\u00d73=447\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop script: replace the date line and every\n# three-digit-by-one-digit multiplication answer in the practice table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2024-12-11 Wednesday\"; New = \"2024-12-12 Thursday\"},\n    @{Old = \"111\u00d76=666\";  New = \"840\u00d79=7560\"},\n    @{Old = \"534\u00d77=3738\"; New = \"832\u00d79=7488\"},\n    @{Old = \"468\u00d73=1404\"; New = \"223\u00d74=892\"},\n    @{Old = \"390\u00d72=780\";  New = \"858\u00d79=7722\"},\n    @{Old = \"891\u00d78=7128\"; New = \"105\u00d79=945\"},\n    @{Old = \"292\u00d79=2628\"; New = \"523\u00d76=3138\"},\n    @{Old = \"973\u00d75=4865\"; New = \"528\u00d79=4752\"},\n    @{Old = \"379\u00d72=758\";  New = \"860\u00d75=4300\"},\n    @{Old = \"380\u00d77=2660\"; New = \"921\u00d73=2763\"},\n    @{Old = \"704\u00d77=4928\"; New = \"485\u00d76=2910\"},\n    @{Old = \"470\u00d74=1880\"; New = \"383\u00d76=2298\"},\n    @{Old = \"545\u00d75=2725\"; New = \"144\u00d74=576\"},\n    @{Old = \"755\u00d76=4530\"; New = \"419\u00d76=2514\"},\n    @{Old = \"389\u00d79=3501\"; New = \"164\u00d75=820\"},\n    @{Old = \"801\u00d78=6408\"; New = \"401\u00d74=1604\"},\n    @{Old = \"337\u00d77=2359\"; New = \"975\u00d74=3900\"},\n    @{Old = \"345\u00d76=2070\"; New = \"248\u00d72=496\"},\n    @{Old = \"486\u00d78=3888\"; New = \"491\u00d74=1964\"},\n    @{Old = \"891\u00d79=8019\"; New = \"125\u00d73=375\"},\n    @{Old = \"446\u00d77=3122\"; New = \"703\u00d76=4218\"},\n    @{Old = \"696\u00d74=2784\"; New = \"242\u00d73=726\"},\n    @{Old = \"205\u00d78=1640\"; New = \"536\u00d77=3752\"},\n    @{Old = \"313\u00d72=626\";  New = \"456\u00d74=1824\"},\n    @{Old = \"380\u00d73=1140\"; New = \"305\u00d73=915\"},\n    @{Old = \"287\u00d72=574\";  New = \"149\u00d73=447\"}\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
